$wb = $excel.ActiveWorkbook

# --- "08-BC" sheet: post the three boundary-condition sections ---
$bc = $wb.Worksheets.Item("08-BC")

# Row 3 already exists with "08-BC.S-03" in A3; rewrite it to the correct
# section id and add its "assigned_region" value in column B.
$bc.Range("A3").Value = "08-BC.S-01"
$bc.Range("B3").Value = "assigned_region"

$bc.Range("A4").Value = "08-BC.S-02"
$bc.Range("B4").Value = "liquid_phase"

$bc.Range("A5").Value = "08-BC.S-03"
$bc.Range("B5").Value = "solid_phase"

# Widen column A to fit the longer section labels (target stored width
# 19.28515625 chars; the host quantizes ColumnWidth to 1/6-character
# steps, so 18.5 is the closest input that lands on the nearest
# representable stored width).
$bc.Columns.Item(1).ColumnWidth = 18.5

# --- "provenance" sheet: refresh the recorded timestamp ---
$prov = $wb.Worksheets.Item("provenance")
$prov.Range("B12").Value = 43435.44090384839
